$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in new row of data: date, event text, author
$ws.Range("B5").Value = (Get-Date -Year 2021 -Month 3 -Day 25).Date
$ws.Range("C5").Value = "Enregistement du score"
$ws.Range("D5").Value = "Arthru Bottemanne"

# Update the active selection to F11
$ws.Range("F11").Select()
